$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 = "Save", matching style of existing header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for column H, rows 2-7
$values = @(0, 0, 0, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
